$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated K column (column G) values - regenerated using K instead of Strike#
$kValues = @{
    2  = 3
    3  = 2
    4  = 4
    5  = 8
    6  = 3
    7  = 0
    8  = 2
    9  = 6
    10 = 1
    11 = 0
    12 = 3
    13 = 3
    14 = 1
    15 = 1
    16 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
